$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Los años"
$ws.Range("B7").Value = "Virginia Woolf"
$ws.Range("C7").Value = "Lumen"
